$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames: replace spaces/punctuation with underscores, drop units ---
$ws.Range("A1").Value = "Mombo_ShotID"
$ws.Range("G1").Value = "Ball_mph"
$ws.Range("H1").Value = "Club_mph"
$ws.Range("I1").Value = "Smash_Factor"
$ws.Range("J1").Value = "Carry_yds"
$ws.Range("K1").Value = "Total_yds"
$ws.Range("L1").Value = "Roll_yds"
$ws.Range("M1").Value = "Swing_H"
$ws.Range("N1").Value = "Spin_rpm"
$ws.Range("O1").Value = "Height_ft"
$ws.Range("P1").Value = "Time_s"
$ws.Range("Q1").Value = "AOA"
$ws.Range("R1").Value = "Spin_Loft"
$ws.Range("S1").Value = "Swing_V"
$ws.Range("T1").Value = "Spin_Axis"
$ws.Range("U1").Value = "Lateral_yds"
$ws.Range("V1").Value = "Shot_Type"
$ws.Range("W1").Value = "FTP"
$ws.Range("X1").Value = "FTT"
$ws.Range("Y1").Value = "Dynamic_Loft"
$ws.Range("Z1").Value = "Club_Path"
$ws.Range("AA1").Value = "Launch_H"
$ws.Range("AB1").Value = "Launch_V"
$ws.Range("AC1").Value = "Low_Point_ftin"
$ws.Range("AD1").Value = "DescentV"
$ws.Range("AE1").Value = "Curve_Dist_yds"
$ws.Range("AF1").Value = "Lateral_Impact_in"
$ws.Range("AG1").Value = "Vertical_Impact_in"
$ws.Range("AJ1").Value = "Unnamed_35"
$ws.Range("AK1").Value = "Unnamed_36"
$ws.Range("AL1").Value = "Unnamed_37"
$ws.Range("AM1").Value = "Unnamed_38"
$ws.Range("AN1").Value = "Unnamed_39"
$ws.Range("AO1").Value = "Unnamed_40"

# --- Data changes for rows 2-6 ---
# Column T (Spin_Axis): numeric value -> text value with " L" suffix
$ws.Range("T2").Value = "13.8 L"
$ws.Range("T3").Value = "9.5 L"
$ws.Range("T4").Value = "7.3 L"
$ws.Range("T5").Value = "12.2 L"
$ws.Range("T6").Value = "3.3 L"

# Column AA (Launch_H): numeric value -> text value with " L" suffix
$ws.Range("AA2").Value = "4.9 L"
$ws.Range("AA3").Value = "4.2 L"
$ws.Range("AA4").Value = "4.8 L"
$ws.Range("AA5").Value = "6.6 L"
$ws.Range("AA6").Value = "4.9 L"

# Column AC (Low_Point_ftin): text value with trailing inch mark -> plain numeric value
$ws.Range("AC2").Value = 4.6
$ws.Range("AC3").Value = 5.9
$ws.Range("AC4").Value = 5.3
$ws.Range("AC5").Value = 4.6
$ws.Range("AC6").Value = 5.2
